# create ram user from execl
# Rewrites the "group" sheet rows: reorders the existing groups and
# appends a new "java_custom," policy entry to the "java" row.

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("group")

$data = @(
    @("china-infra", "AdministratorAccess,"),
    @("group-go", "AliyunLogReadOnlyAccess,AliyunCloudMonitorReadOnlyAccess,AliyunDysmsReadOnlyAccess,ContainerRegistryPullPushAccess,StagingSwarmAccess,GoSquadCdnPolicy,GoSquadOssPolicy,"),
    @("ruby", "AliyunRDSReadOnlyAccess,AliyunLogReadOnlyAccess,AliyunCSReadOnlyAccess,AliyunContainerRegistryFullAccess,StagingRDSFullAccess,"),
    @("java", "AliyunOSSReadOnlyAccess,AliyunLogReadOnlyAccess,AliyunCSReadOnlyAccess,AliyunContainerRegistryFullAccess,AliyunMQReadOnlyAccess,java_custom,"),
    @("billing-squad", "StagingSwarmAccess,ProductionSwarmAccess,BillingSquadOssPolicy,BillingSquadCdnPolicy,ChinaSpacestationOssPolicy,"),
    @("engagement", "AliyunMQPubOnlyAccess,AliyunMQSubOnlyAccess,AliyunContainerRegistryReadOnlyAccess,StagingEngagementOSSAccess,ContainerRegistryPullPushAccess,StagingSwarmAccess,"),
    @("group-member-to-member", "AliyunLogReadOnlyAccess,AliyunCloudMonitorReadOnlyAccess,AliyunDysmsReadOnlyAccess,ContainerRegistryPullPushAccess,StagingSwarmAccess,M2MSquadOssPolicy,M2MSquadCdnPolicy,"),
    @("growth-and-billing", "AliyunContainerRegistryReadOnlyAccess,StagingWWWOSSAccess,StagingRDSFullAccess,ContainerRegistryPullPushAccess,ProductionWWWOSSAccess,StagingSwarmAccess,"),
    @("growth-and-billing-admin", "AliyunOSSFullAccess,AliyunRDSFullAccess,AliyunKvstoreFullAccess,AliyunCloudMonitorFullAccess,AliyunDysmsFullAccess,AliyunContainerRegistryFullAccess,ListRAMUsers,ProductionSwarmAccess,"),
    @("bastion-operation", "AliyunYundunBastionHostReadOnlyAccess,AliyunYundunBastionHostOperateOnlyAccess,"),
    @("china-data-squad", "AliyunRDSReadOnlyAccess,AliyunDataWorksFullAccess,"),
    @("group-member-to-space", "AliyunLogReadOnlyAccess,AliyunCloudMonitorReadOnlyAccess,StagingSwarmAccess,M2SSquadOssPolicy,M2SSquadCdnPolicy,"),
    @("growth-squad", "AliyunLogReadOnlyAccess,AliyunCloudMonitorReadOnlyAccess,AliyunDysmsReadOnlyAccess,ContainerRegistryPullPushAccess,StagingSwarmAccess,GrowthSquadOssPolicy,GrowthSquadCdnPolicy,"),
    @("recommendation-squad", "ContainerRegistryPullPushAccess,"),
    @("china-spacestation-developers", "ChinaSpacestationOssPolicy,"),
    @("dev-platform", "AdministratorAccess,ContainerRegistryPullPushAccess,"),
    @("engagement-admin", "StagingECSROAccess,ProductionECSROAccess,StagingSwarmAccess,ProductionSwarmAccess,")
)

for ($i = 0; $i -lt $data.Length; $i++) {
    $row = $i + 1
    $ws.Cells.Item($row, 1).Value = $data[$i][0]
    $ws.Cells.Item($row, 2).Value = $data[$i][1]
}
